$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.488.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.030.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.022.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.445"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000222"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.518.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.349.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.020.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.108"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.672"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "460.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.154.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0801"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0387"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.26%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.40%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.244"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.109"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0500"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.78%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.65%  "
